# Update res_bus/vm_pu values for the "380 kV" case (Case_2_20)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (in sheet order, column H / row index 8 is intentionally absent from the data)
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# New vm_pu values per data row (rows 2-25), one array per row aligned with $cols
$newValues = @{
    2 = @(1.02, 1.026573563065708, 1.034124669515817, 0.9926147277508489, 1.041235863140376, 1.031049196888618, 1.031736249692989, 1.036925044196291, 0.9955398523336033, 1.044015950788455, 1.033201433509366)
    3 = @(1.02, 1.027495907404926, 1.034815832032269, 0.9936372048519304, 1.042121140327257, 1.031170432803451, 1.032298551266623, 1.037425644031639, 0.9963617723202692, 1.044711618174526, 1.03376453361572)
    4 = @(1.02, 1.028093147393348, 1.035263256265587, 0.9942998659930995, 1.042694665044415, 1.031247571963841, 1.032662210493352, 1.037749085541785, 0.9968940712668345, 1.045161816684257, 1.034128709280264)
    5 = @(1.02, 1.028344326920691, 1.035451398788775, 0.9945786998346017, 1.042935938748832, 1.031279687421343, 1.032815046916827, 1.037884944301923, 0.997117960005301, 1.045351092014607, 1.034281762748991)
    6 = @(1.02, 1.028386506891759, 1.035482991363012, 0.9946255319796338, 1.042976459232158, 1.031285061327415, 1.032840706102407, 1.037907748747595, 0.9971555583673453, 1.045382872833029, 1.034307458373557)
    7 = @(1.02, 1.02809650327554, 1.035265770058805, 0.9943035907982488, 1.04269788831335, 1.031248002326064, 1.032664252882831, 1.037750901350601, 0.9968970624462087, 1.045164345745496, 1.03413075457017)
    8 = @(1.02, 1.02688518572489, 1.034358209657635, 0.9929600610674301, 1.041534902661039, 1.031090439742833, 1.031926320344751, 1.037094322977899, 0.995817528259106, 1.044251042607031, 1.033391774083244)
    9 = @(1.02, 1.024753965381754, 1.032760539237387, 0.9906006454969559, 1.039490935660725, 1.030802797228527, 1.030624598370008, 1.035933720892148, 0.9939188001724441, 1.042642167453311, 1.032088203513987)
    10 = @(1.02, 1.023335417299284, 1.031696568793661, 0.989033133672735, 1.038131988231047, 1.030604347639604, 1.029755909357289, 1.03515761381342, 0.9926553831429383, 1.041569986716987, 1.03121828086323)
    11 = @(1.02, 1.0227217205313, 1.031236147740429, 0.988355674866747, 1.037544445094228, 1.030516837485123, 1.029379562600141, 1.034821001999965, 0.9921088820399291, 1.04110583256173, 1.030841399650519)
    12 = @(1.02, 1.022493848931058, 1.031065171198582, 0.9881042295826724, 1.037326340658, 1.030484095451401, 1.029239741790683, 1.034695887489517, 0.9919059725120875, 1.040933442441586, 1.030701380279474)
    13 = @(1.02, 1.022542724404405, 1.031101844216873, 0.9881581567098651, 1.037373118655376, 1.030491129438834, 1.029269735119459, 1.034722728665232, 0.9919494934313052, 1.040970419899359, 1.030731416202217)
    14 = @(1.02, 1.022702882900872, 1.03122201385304, 0.9883348863814464, 1.037526413745728, 1.030514135847287, 1.029368005551607, 1.034810661659865, 0.9920921077337197, 1.041091582378174, 1.030829826189651)
    15 = @(1.02, 1.022801572825163, 1.031296060251784, 0.9884438009545853, 1.037620881877099, 1.030528279483882, 1.029428549412427, 1.034864829260551, 0.9921799884222134, 1.041166236924647, 1.030890456029698)
    16 = @(1.02, 1.023376157876436, 1.031727131558394, 0.9890781214508737, 1.038171000416688, 1.03061012215376, 1.029780882149971, 1.035179942086664, 0.9926916645766087, 1.041600793468487, 1.031243289120142)
    17 = @(1.02, 1.02373672621567, 1.031997608603892, 0.989476357848556, 1.038516314614056, 1.030661037220205, 1.030001838668637, 1.035377456930372, 0.9930127773699352, 1.041873408918345, 1.031464559422408)
    18 = @(1.02, 1.023947092040888, 1.032155400703144, 0.9897087662937556, 1.038717816319698, 1.030690582617222, 1.030130699631832, 1.035492610632731, 0.9932001317071769, 1.042032431152201, 1.031593603382952)
    19 = @(1.02, 1.024018830211157, 1.032209208335798, 0.9897880325774034, 1.03878653769692, 1.030700630951134, 1.030174634597584, 1.035531865984849, 0.9932640239640975, 1.04208665532427, 1.031637600741395)
    20 = @(1.02, 1.023698035216734, 1.031968586130239, 0.9894336180360679, 1.038479256785409, 1.030655590285317, 1.029978134107831, 1.035356270968764, 0.9929783193494215, 1.041844158790549, 1.031440821198408)
    21 = @(1.02, 1.022655717924516, 1.031186625654418, 0.9882828385668249, 1.037481268416745, 1.03050736756776, 1.029339068147643, 1.034784769836614, 0.9920501090198102, 1.041055902557329, 1.030800847691254)
    22 = @(1.02, 1.022000850666293, 1.030695233328598, 0.9875604150241495, 1.036854576759437, 1.030412803895215, 1.028937095484971, 1.034424971107963, 0.9914670000341481, 1.040560395888168, 1.030398304181294)
    23 = @(1.02, 1.022347962456666, 1.030955704814373, 0.9879432794643023, 1.037186723113984, 1.030463063583734, 1.029150204255886, 1.034615751721222, 0.991776070289318, 1.040823063312703, 1.030611715591107)
    24 = @(1.02, 1.023715517851833, 1.03198170005228, 0.9894529299347244, 1.03849600135816, 1.030658051991984, 1.029988845238033, 1.035365844157869, 0.9929938892766442, 1.041857375632687, 1.031451547539643)
    25 = @(1.02, 1.025304542237465, 1.033173380398368, 0.9912096547607049, 1.040018705188613, 1.03087834096378, 1.030961283249053, 1.036234187345478, 0.9944092447426414, 1.043058035037453, 1.032425366524179)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value2 = $vals[$i]
    }
}
